# Generate Report for Handoff
# Replace the old GUID-named handoff file references (167f055f-...) with the
# newly generated ones (d88e5db0-...), and bump the associated timestamps
# to reflect the new handoff/xliff-generation run.

$wb = $excel.ActiveWorkbook

$oldGuid = "167f055f-b88f-4b9a-9086-6b3f7126aadd"
$newGuid = "d88e5db0-5aa9-468e-87d8-9526b2ea6769"

$oldHash = "090a360ebddec72f55a84902313e81e782f66f2f"
$newHash = "2e5cee938237a608871ae136aecc6501d2ac445d"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$newDisplayMd = "e2e\$newGuid.md"
$wsOverview.Range("B2").Value = $newDisplayMd
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newDisplayMd
}

$wsOverview.Range("G2").Value = "2016-08-26 00:57:20"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$newFileNameMd = "$newGuid.md"
$wsZhCn.Range("A2").Value = $newFileNameMd
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $newFileNameMd
}

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 00:57:16"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFileNameMd
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $newFileNameMd
}

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 00:57:20"
